$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "24.313.34"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +0.76%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.668.06"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +1.67%  "

$ws.Range("E4").Value = "  +0.35%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "311.98"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +1.52%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "1.003"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +0.27%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.3957"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +1.85%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3923"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +1.61%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "52.21"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +5.34%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "1.387"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +3.01%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "1.003"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +0.40%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.08568"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -1.63%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "24.38"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +3.27%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "7.290"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +2.79%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "8.005"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +7.46%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.00001332"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +3.42%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "1.666.90"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +1.87%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "95.33"
$cell.Style = "Normal"

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.07033"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +2.18%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "20.50"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.42%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "6.992"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +1.52%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.26%  "

$ws.Range("E23").Value = "  +1.33%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "24.325.77"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.85%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.532"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +9.04%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "3.101"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +14.46%  "

$ws.Range("E27").Value = "  +1.01%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "156.96"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -0.32%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "142.68"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +1.79%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "5.340"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -0.13%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "8.019"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -6.44%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "2.548"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +5.68%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "1.846.89"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +1.68%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.059"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +12.22%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.03099"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +7.59%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.08212"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +2.80%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "6.874"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -0.34%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "11.14"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +12.65%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.2757"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +3.43%  "

$ws.Range("E40").Value = "  +1.06%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "13.82"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +6.33%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.7674"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +1.86%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "1.438"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -0.94%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "16.64"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +4.71%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.7073"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +2.87%  "

$ws.Range("E46").Value = "  +2.81%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "4.122"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +0.99%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +0.25%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.08416"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.47%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "136.30"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +2.77%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "1.259"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +0.15%  "
